# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2410"
#   "<header>_new" -> "<header>_FV2504"
# and turn the data range into a formatted/filterable Excel Table with a
# frozen header row, mirroring how the exporter now renders these AHB
# diff sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the header cells (row 1) to use the format-version suffixes
#    instead of the generic "_old" / "_new" suffixes.
# ---------------------------------------------------------------------
$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $text = $cell.Value()
    $cell.Value = $text.Replace("_old", "_FV2410")
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $text = $cell.Value()
    $cell.Value = $text.Replace("_new", "_FV2504")
}

# ---------------------------------------------------------------------
# 2. Convert the used range into an Excel Table ("Table1") including an
#    autofilter, matching the exporter's new xlsx output.
# ---------------------------------------------------------------------
$dataRange = $ws.Range("A1:U61")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# ---------------------------------------------------------------------
# 3. Freeze the header row so it stays visible while scrolling.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

Write-Output "Header renaming, table creation and freeze panes applied."
